$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric need to be forced to text so Excel
# does not auto-convert them (e.g. "242.20" -> 242.2, "0.06610" -> 0.0661).
# Set Text number format, assign the literal string, then clear the format
# again so the cell keeps its original (default) style, matching the source diff.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D16", "D18", "D19", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.038.45"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.878.58"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "242.20"
$ws.Range("E5").Value = "  -2.51%  "
$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "0.4956"
$ws.Range("E7").Value = "  -0.95%  "
$ws.Range("D8").Value = "44.48"
$ws.Range("E8").Value = "  -2.74%  "
$ws.Range("D9").Value = "0.2922"
$ws.Range("E9").Value = "  +2.27%  "
$ws.Range("D10").Value = "0.06610"
$ws.Range("E10").Value = "  +1.06%  "
$ws.Range("D11").Value = "1.879.52"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "16.74"
$ws.Range("E12").Value = "  -2.06%  "
$ws.Range("D13").Value = "0.07170"
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("D14").Value = "0.6624"
$ws.Range("E14").Value = "  -0.59%  "
$ws.Range("D15").Value = "86.04"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").Value = "4.841"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").Value = "30.001.33"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "0.000007853"
$ws.Range("E18").Value = "  +4.46%  "
$ws.Range("D19").Value = "0.9997"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").Value = "2.121.08"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "4.755"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("D25").Value = "9.100"
$ws.Range("E25").Value = "  +0.94%  "
$ws.Range("D26").Value = "150.53"
$ws.Range("E26").Value = "  +3.50%  "
$ws.Range("D27").Value = "134.45"
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("D28").Value = "16.72"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "1.908"
$ws.Range("E29").Value = "  -2.18%  "
$ws.Range("D30").Value = "1.371"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Value = "4.163"
$ws.Range("E31").Value = "  -0.77%  "
$ws.Range("D32").Value = "0.08690"
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("D33").Value = "3.943"
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("E34").Value = "  -1.37%  "
$ws.Range("E35").Value = "  +2.06%  "
$ws.Range("D36").Value = "1.099"
$ws.Range("E36").Value = "  -2.91%  "
$ws.Range("D37").Value = "2.654"
$ws.Range("E37").Value = "  -1.00%  "
$ws.Range("D38").Value = "2.697"
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("D39").Value = "2.180"
$ws.Range("E39").Value = "  -5.10%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "0.9326"
$ws.Range("E40").Value = "  -2.86%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.01685"
$ws.Range("E41").Value = "  +3.15%  "
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("D43").Value = "0.9991"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").Value = "0.4186"
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("D45").Value = "101.72"
$ws.Range("E45").Value = "  -2.78%  "
$ws.Range("D46").Value = "7.410"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("D47").Value = "0.1257"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").Value = "0.05668"
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("D49").Value = "32.45"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "55.80"
$ws.Range("E50").Value = "  +1.77%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.335"
$ws.Range("E51").Value = "  +0.16%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
